$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns - copy formatting (bold, border,
# centered alignment) from the existing header cell H1, then set values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for new columns I (I0) and J (IF)
$data = @(
    @(2, 9, 9),
    @(3, 10, 10),
    @(4, 7, 7),
    @(5, 4, 5),
    @(6, 5, 6),
    @(7, 6, 7),
    @(8, 1, 1),
    @(9, 6, 6),
    @(10, 9, 9),
    @(11, 8, 9),
    @(12, 8, 8),
    @(13, 8, 8),
    @(14, 6, 7),
    @(15, 9, 9),
    @(16, 9, 9),
    @(17, 8, 9),
    @(18, 7, 8),
    @(19, 8, 8),
    @(20, 6, 6)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
